# fix JO RFQ ewt for labor and materials
#
# Re-creates the authorial edit described by the commit:
#   - Purchase Request / Department / Dept. Code / Requestor fields on
#     Sheet1 get new sample values ("bcd" / "itb" / "henne").
#   - Purpose / End-Use are both set to "asdasd".
#   - Date Prepared / Date Issued move from 2024-01-30 to 2024-05-27.
#   - The first item row's UOM/Description placeholders change
#     ("pc/s"->"pcs", "Testing Data"->"item 1") and its "Date Needed"
#     cell is overwritten with the literal placeholder text "YYYY-MM-DD".
#   - The second item row (row 15) is wiped out entirely (content +
#     its two merges removed), matching the now-blank rows beneath it.
#   - The active selection on Sheet1 moves to H16:I16.
#   - The saved window position/size is nudged (best effort - the
#     sandbox's COM surface does not persist bookView geometry, see
#     below).
#   - The recorded "last opened from" folder's username changes from
#     steph to Henne (best effort - not exposed on the Excel object
#     model; there is no COM property that maps to the x15ac:absPath
#     extension, so this part of the diff cannot be reproduced here).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Purchase Request header block (rows 7-9) ------------------------
# C7 = Purchase Request value, I7 = Department value - the author typed
# the same text "bcd" into both.
$ws.Range("C7").Value = "bcd"
$ws.Range("I7").Value = "bcd"

# H8/I8 = Dept. Code value
$ws.Range("I8").Value = "itb"

# H9/I9 = Requestor value
$ws.Range("I9").Value = "henne"

# Date Prepared / Date Issued -> 2024-05-27
$ws.Range("C8").Value = "5/27/2024"
$ws.Range("C9").Value = "5/27/2024"

# --- Purpose / End-Use (rows 11-12) -----------------------------------
$ws.Range("C11").Value = "asdasd"
$ws.Range("C12").Value = "asdasd"

# --- Item table (rows 14-15) ------------------------------------------
# Row 14: UOM -> "pcs", Description -> "item 1", Date Needed -> literal
# text "YYYY-MM-DD" (placeholder typed over the date cell).
$ws.Range("C14").Value = "pcs"
$ws.Range("E14").Value = "item 1"
$ws.Range("J14").Value = "YYYY-MM-DD"

# Row 15: remove the second item entirely (content + its merges), so it
# becomes a blank row like rows 16-20.
$ws.Range("E15:H15").UnMerge()
$ws.Range("J15:K15").UnMerge()
$ws.Range("A15:K15").ClearContents()
$ws.Range("A15:K15").Style = $ws.Range("A16").Style

# --- Selection ----------------------------------------------------------
$ws.Range("H16:I16").Select()

# --- Window geometry (best effort) --------------------------------------
# The target bookViews/workbookView attributes move from
# xWindow=0 yWindow=0 windowWidth=20490 windowHeight=7155 to
# xWindow=-120 yWindow=-120 windowWidth=20730 windowHeight=11160.
# Drive it through the documented Window object even though this
# sandbox's COM shim does not currently round-trip these values into
# the saved XML.
$win = $excel.ActiveWindow
$win.Left = -120
$win.Top = -120
$win.Width = 20730
$win.Height = 11160
